# Append new log_activity rows (394-434) captured after the date/timezone
# fix, mirroring the commit "update: try to fix date or timezone api".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(394, 1).Value = 'I3KNC0'
$ws.Cells.Item(394, 2).Value = '2024-11-22 21:06:24'
$ws.Cells.Item(394, 3).Value = 'GET /'
$ws.Cells.Item(394, 4).Value = 401
$ws.Cells.Item(394, 5).Value = $false
$ws.Cells.Item(394, 6).Value = 'Eitss... mau ngapain? Akses terbatas!'

$ws.Cells.Item(395, 1).Value = 'MLUJ5S'
$ws.Cells.Item(395, 2).Value = '2024-11-22 21:06:24'
$ws.Cells.Item(395, 3).Value = 'GET /favicon.ico'
$ws.Cells.Item(395, 4).Value = 404
$ws.Cells.Item(395, 5).Value = $false
$ws.Cells.Item(395, 6).Value = 'Not Found'

$ws.Cells.Item(396, 1).Value = 'BBIG20'
$ws.Cells.Item(396, 2).Value = '2024-11-22 21:06:37'
$ws.Cells.Item(396, 3).Value = 'GET /'
$ws.Cells.Item(396, 4).Value = 200
$ws.Cells.Item(396, 5).Value = $true
$ws.Cells.Item(396, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:06:37'

$ws.Cells.Item(397, 1).Value = 'B4IT6K'
$ws.Cells.Item(397, 2).Value = '2024-11-22 19:06:53'
$ws.Cells.Item(397, 3).Value = 'GET /'
$ws.Cells.Item(397, 4).Value = 200
$ws.Cells.Item(397, 5).Value = $true
$ws.Cells.Item(397, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:06:53'

$ws.Cells.Item(398, 1).Value = 'PD78F5'
$ws.Cells.Item(398, 2).Value = '2024-11-22 19:06:56'
$ws.Cells.Item(398, 3).Value = 'GET /'
$ws.Cells.Item(398, 4).Value = 200
$ws.Cells.Item(398, 5).Value = $true
$ws.Cells.Item(398, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:06:56'

$ws.Cells.Item(399, 1).Value = 'LBG84X'
$ws.Cells.Item(399, 2).Value = '2024-11-22 19:10:18'
$ws.Cells.Item(399, 3).Value = 'GET /'
$ws.Cells.Item(399, 4).Value = 200
$ws.Cells.Item(399, 5).Value = $true
$ws.Cells.Item(399, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:10:18'

$ws.Cells.Item(400, 1).Value = 'HV2TCV'
$ws.Cells.Item(400, 2).Value = '2024-11-22 19:10:19'
$ws.Cells.Item(400, 3).Value = 'GET /'
$ws.Cells.Item(400, 4).Value = 200
$ws.Cells.Item(400, 5).Value = $true
$ws.Cells.Item(400, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:10:19'

$ws.Cells.Item(401, 1).Value = 'NLFWZ4'
$ws.Cells.Item(401, 2).Value = '2024-11-22 19:10:19'
$ws.Cells.Item(401, 3).Value = 'GET /'
$ws.Cells.Item(401, 4).Value = 200
$ws.Cells.Item(401, 5).Value = $true
$ws.Cells.Item(401, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:10:19'

$ws.Cells.Item(402, 1).Value = 'YBRISH'
$ws.Cells.Item(402, 2).Value = '2024-11-22 19:10:20'
$ws.Cells.Item(402, 3).Value = 'GET /'
$ws.Cells.Item(402, 4).Value = 200
$ws.Cells.Item(402, 5).Value = $true
$ws.Cells.Item(402, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:10:20'

$ws.Cells.Item(403, 1).Value = 'OY0BNU'
$ws.Cells.Item(403, 2).Value = '2024-11-22 19:10:20'
$ws.Cells.Item(403, 3).Value = 'GET /'
$ws.Cells.Item(403, 4).Value = 200
$ws.Cells.Item(403, 5).Value = $true
$ws.Cells.Item(403, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:10:20'

$ws.Cells.Item(404, 1).Value = '6M3GHU'
$ws.Cells.Item(404, 2).Value = '2024-11-22 19:10:21'
$ws.Cells.Item(404, 3).Value = 'GET /'
$ws.Cells.Item(404, 4).Value = 200
$ws.Cells.Item(404, 5).Value = $true
$ws.Cells.Item(404, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:10:21'

$ws.Cells.Item(405, 1).Value = '8ZEA8Z'
$ws.Cells.Item(405, 2).Value = '2024-11-22 19:10:21'
$ws.Cells.Item(405, 3).Value = 'GET /'
$ws.Cells.Item(405, 4).Value = 200
$ws.Cells.Item(405, 5).Value = $true
$ws.Cells.Item(405, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:10:21'

$ws.Cells.Item(406, 1).Value = '19YEPX'
$ws.Cells.Item(406, 2).Value = '2024-11-22 19:10:21'
$ws.Cells.Item(406, 3).Value = 'GET /'
$ws.Cells.Item(406, 4).Value = 200
$ws.Cells.Item(406, 5).Value = $true
$ws.Cells.Item(406, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:10:21'

$ws.Cells.Item(407, 1).Value = '3X854H'
$ws.Cells.Item(407, 2).Value = '2024-11-22 19:10:29'
$ws.Cells.Item(407, 3).Value = 'GET /'
$ws.Cells.Item(407, 4).Value = 200
$ws.Cells.Item(407, 5).Value = $true
$ws.Cells.Item(407, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:10:29'

$ws.Cells.Item(408, 1).Value = '9CGKO1'
$ws.Cells.Item(408, 2).Value = '2024-11-22 19:10:30'
$ws.Cells.Item(408, 3).Value = 'GET /'
$ws.Cells.Item(408, 4).Value = 200
$ws.Cells.Item(408, 5).Value = $true
$ws.Cells.Item(408, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:10:30'

$ws.Cells.Item(409, 1).Value = 'RJ4P7W'
$ws.Cells.Item(409, 2).Value = '2024-11-22 19:20:10'
$ws.Cells.Item(409, 3).Value = 'GET /'
$ws.Cells.Item(409, 4).Value = 200
$ws.Cells.Item(409, 5).Value = $true
$ws.Cells.Item(409, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:20:10'

$ws.Cells.Item(410, 1).Value = 'KVITQC'
$ws.Cells.Item(410, 2).Value = '2024-11-22 19:24:32'
$ws.Cells.Item(410, 3).Value = 'GET /'
$ws.Cells.Item(410, 4).Value = 200
$ws.Cells.Item(410, 5).Value = $true
$ws.Cells.Item(410, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:24:32'

$ws.Cells.Item(411, 1).Value = 'CX6JB8'
$ws.Cells.Item(411, 2).Value = '2024-11-22 19:25:25'
$ws.Cells.Item(411, 3).Value = 'GET /'
$ws.Cells.Item(411, 4).Value = 200
$ws.Cells.Item(411, 5).Value = $true
$ws.Cells.Item(411, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:25:25'

$ws.Cells.Item(412, 1).Value = 'C85Y1Q'
$ws.Cells.Item(412, 2).Value = '2024-11-22 19:25:33'
$ws.Cells.Item(412, 3).Value = 'GET /'
$ws.Cells.Item(412, 4).Value = 200
$ws.Cells.Item(412, 5).Value = $true
$ws.Cells.Item(412, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:25:33'

$ws.Cells.Item(413, 1).Value = '0OIR6R'
$ws.Cells.Item(413, 2).Value = '2024-11-22 19:25:34'
$ws.Cells.Item(413, 3).Value = 'GET /'
$ws.Cells.Item(413, 4).Value = 200
$ws.Cells.Item(413, 5).Value = $true
$ws.Cells.Item(413, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:25:34'

$ws.Cells.Item(414, 1).Value = 'UECCLM'
$ws.Cells.Item(414, 2).Value = '2024-11-22 19:26:43'
$ws.Cells.Item(414, 3).Value = 'GET /'
$ws.Cells.Item(414, 4).Value = 200
$ws.Cells.Item(414, 5).Value = $true
$ws.Cells.Item(414, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:26:43'

$ws.Cells.Item(415, 1).Value = 'R9YFUY'
$ws.Cells.Item(415, 2).Value = '2024-11-22 19:26:51'
$ws.Cells.Item(415, 3).Value = 'GET /'
$ws.Cells.Item(415, 4).Value = 200
$ws.Cells.Item(415, 5).Value = $true
$ws.Cells.Item(415, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:26:51'

$ws.Cells.Item(416, 1).Value = '4P1D98'
$ws.Cells.Item(416, 2).Value = '2024-11-22 19:27:21'
$ws.Cells.Item(416, 3).Value = 'GET /datasets/list'
$ws.Cells.Item(416, 4).Value = 200
$ws.Cells.Item(416, 5).Value = $true
$ws.Cells.Item(416, 6).Value = 'Datasets ditemukan.
###
Datasets:[''Beasiswa.pdf'', ''Daya Tampung Mahasiswa Baru.pdf'', ''FAQ DOSEN.pdf'', ''FAQ MAHASISWA.pdf'', ''FAQ PENERIMAAN MAHASISWA SMBJM.pdf'', ''FAQ REMUNERASI.pdf'', ''FAQ UMUM.pdf'', ''Fasilitas.pdf'', ''FBS.pdf'', ''FE.pdf'', ''FHIS.pdf'', ''FIP.pdf'', ''FK.pdf'', ''FMIPA.pdf'', ''FOK.pdf'', ''FTK.pdf'', ''Jadwal PMB.pdf'', ''Kurikulum.pdf'', ''PANDUAN PENDAFTARAN KEMBALI SMBJM dengan skor UTBK-SNBT dan Prestasi-signed.pdf'', ''Panduan-Pendaftaran-Kembali-SMBJM-CBT_2024rv1.pdf'', ''Panduan-Pendaftaran-Kembali-SNBP-2024-v4-wa-signed.pdf'', ''Panduan-Pendaftaran-Kembali-SNBT-2024_rev_signed.pdf'', ''Pascasarjana.pdf'', ''Pedoman-Studi-2017.pdf'', ''Prosedur Pembayaran UKT.pdf'', ''SMBJM.pdf'', ''SNBP.pdf'', ''SNBT.pdf'', ''UKM.pdf'', ''Umum.pdf'', ''UPA TIK.pdf'']'

$ws.Cells.Item(417, 1).Value = '0ZBJ7W'
$ws.Cells.Item(417, 2).Value = '2024-11-22 21:27:54'
$ws.Cells.Item(417, 3).Value = 'GET /'
$ws.Cells.Item(417, 4).Value = 200
$ws.Cells.Item(417, 5).Value = $true
$ws.Cells.Item(417, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:27:54'

$ws.Cells.Item(418, 1).Value = 'T7R75F'
$ws.Cells.Item(418, 2).Value = '2024-11-22 21:29:11'
$ws.Cells.Item(418, 3).Value = 'GET /'
$ws.Cells.Item(418, 4).Value = 200
$ws.Cells.Item(418, 5).Value = $true
$ws.Cells.Item(418, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:29:11'

$ws.Cells.Item(419, 1).Value = '8WRJI1'
$ws.Cells.Item(419, 2).Value = '2024-11-22 21:29:14'
$ws.Cells.Item(419, 3).Value = 'GET /'
$ws.Cells.Item(419, 4).Value = 200
$ws.Cells.Item(419, 5).Value = $true
$ws.Cells.Item(419, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:29:14'

$ws.Cells.Item(420, 1).Value = 'QFG9KH'
$ws.Cells.Item(420, 2).Value = '2024-11-22 21:30:10'
$ws.Cells.Item(420, 3).Value = 'GET /'
$ws.Cells.Item(420, 4).Value = 200
$ws.Cells.Item(420, 5).Value = $true
$ws.Cells.Item(420, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:30:10'

$ws.Cells.Item(421, 1).Value = 'ABY6TG'
$ws.Cells.Item(421, 2).Value = '2024-11-22 21:30:12'
$ws.Cells.Item(421, 3).Value = 'GET /'
$ws.Cells.Item(421, 4).Value = 200
$ws.Cells.Item(421, 5).Value = $true
$ws.Cells.Item(421, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:30:12'

$ws.Cells.Item(422, 1).Value = '4UEU7C'
$ws.Cells.Item(422, 2).Value = '2024-11-22 21:30:12'
$ws.Cells.Item(422, 3).Value = 'GET /'
$ws.Cells.Item(422, 4).Value = 200
$ws.Cells.Item(422, 5).Value = $true
$ws.Cells.Item(422, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:30:12'

$ws.Cells.Item(423, 1).Value = 'O92CB0'
$ws.Cells.Item(423, 2).Value = '2024-11-22 21:30:13'
$ws.Cells.Item(423, 3).Value = 'GET /'
$ws.Cells.Item(423, 4).Value = 200
$ws.Cells.Item(423, 5).Value = $true
$ws.Cells.Item(423, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:30:13'

$ws.Cells.Item(424, 1).Value = 'F5PRKE'
$ws.Cells.Item(424, 2).Value = '2024-11-22 21:30:13'
$ws.Cells.Item(424, 3).Value = 'GET /'
$ws.Cells.Item(424, 4).Value = 200
$ws.Cells.Item(424, 5).Value = $true
$ws.Cells.Item(424, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:30:13'

$ws.Cells.Item(425, 1).Value = 'ZA20DR'
$ws.Cells.Item(425, 2).Value = '2024-11-22 21:30:14'
$ws.Cells.Item(425, 3).Value = 'GET /'
$ws.Cells.Item(425, 4).Value = 200
$ws.Cells.Item(425, 5).Value = $true
$ws.Cells.Item(425, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:30:14'

$ws.Cells.Item(426, 1).Value = 'RB2EC4'
$ws.Cells.Item(426, 2).Value = '2024-11-22 21:30:49'
$ws.Cells.Item(426, 3).Value = 'GET /'
$ws.Cells.Item(426, 4).Value = 200
$ws.Cells.Item(426, 5).Value = $true
$ws.Cells.Item(426, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:30:49'

$ws.Cells.Item(427, 1).Value = 'SXPMBM'
$ws.Cells.Item(427, 2).Value = '2024-11-22 21:30:50'
$ws.Cells.Item(427, 3).Value = 'GET /'
$ws.Cells.Item(427, 4).Value = 200
$ws.Cells.Item(427, 5).Value = $true
$ws.Cells.Item(427, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:30:50'

$ws.Cells.Item(428, 1).Value = 'UMF3VX'
$ws.Cells.Item(428, 2).Value = '2024-11-22 21:32:47'
$ws.Cells.Item(428, 3).Value = 'GET /'
$ws.Cells.Item(428, 4).Value = 401
$ws.Cells.Item(428, 5).Value = $false
$ws.Cells.Item(428, 6).Value = 'Eitss... mau ngapain? Akses terbatas!'

$ws.Cells.Item(429, 1).Value = 'MZYR0S'
$ws.Cells.Item(429, 2).Value = '2024-11-22 21:32:47'
$ws.Cells.Item(429, 3).Value = 'GET /favicon.ico'
$ws.Cells.Item(429, 4).Value = 404
$ws.Cells.Item(429, 5).Value = $false
$ws.Cells.Item(429, 6).Value = 'Not Found'

$ws.Cells.Item(430, 1).Value = 'D6GIIX'
$ws.Cells.Item(430, 2).Value = '2024-11-22 21:32:48'
$ws.Cells.Item(430, 3).Value = 'GET /'
$ws.Cells.Item(430, 4).Value = 401
$ws.Cells.Item(430, 5).Value = $false
$ws.Cells.Item(430, 6).Value = 'Eitss... mau ngapain? Akses terbatas!'

$ws.Cells.Item(431, 1).Value = 'KE2DNM'
$ws.Cells.Item(431, 2).Value = '2024-11-22 21:33:00'
$ws.Cells.Item(431, 3).Value = 'GET /'
$ws.Cells.Item(431, 4).Value = 200
$ws.Cells.Item(431, 5).Value = $true
$ws.Cells.Item(431, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:33:00'

$ws.Cells.Item(432, 1).Value = 'Q5CLLD'
$ws.Cells.Item(432, 2).Value = '2024-11-22 21:36:33'
$ws.Cells.Item(432, 3).Value = 'GET /'
$ws.Cells.Item(432, 4).Value = 401
$ws.Cells.Item(432, 5).Value = $false
$ws.Cells.Item(432, 6).Value = 'Eitss... mau ngapain? Akses terbatas!'

$ws.Cells.Item(433, 1).Value = 'K3UBQV'
$ws.Cells.Item(433, 2).Value = '2024-11-22 21:36:51'
$ws.Cells.Item(433, 3).Value = 'GET /'
$ws.Cells.Item(433, 4).Value = 200
$ws.Cells.Item(433, 5).Value = $true
$ws.Cells.Item(433, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:36:51'

$ws.Cells.Item(434, 1).Value = 'C3LMV9'
$ws.Cells.Item(434, 2).Value = '2024-11-22 21:37:35'
$ws.Cells.Item(434, 3).Value = 'GET /'
$ws.Cells.Item(434, 4).Value = 200
$ws.Cells.Item(434, 5).Value = $true
$ws.Cells.Item(434, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 21:37:35'
